$d = $word.ActiveDocument

# 0) Remove the old "_GoBack" bookmark that currently sits between "Ngay"
#    and "Hoc" - it is about to move to the code-snippet paragraph below.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Locate the target run ("Microsoft.AspNetCore.Identity.EntityFrameworkCore;").
$full = $d.Content
$full.Find.Execute(
    "Microsoft.AspNetCore.Identity.EntityFrameworkCore;", $true, $false, $false,
    $false, $false, $true, 1, $false, "", 0) | Out-Null

# 1) Wrap the trailing ";" in a (temporary, non-collapsed) "_GoBack" bookmark,
#    then delete that character. Deleting text spanned by a bookmark leaves
#    the bookmark collapsed at that position, which is exactly how Word
#    marks the location of the last edit.
$semi = $d.Range($full.End - 1, $full.End)
$semi.Bookmarks.Add("_GoBack") | Out-Null
$semi2 = $d.Range($full.End - 1, $full.End)
$semi2.Text = ""

# 2) Split the run "Microsoft.AspNetCore.Identity.EntityFrameworkCore" into
#    two runs ("Microsoft.AspNetCo" + "re.Identity.EntityFrameworkCore")
#    while keeping the formatting identical. Nudge a character property away
#    and then back to its original value on the second half of the text so
#    Word materializes a run boundary there.
$tail = $d.Content
$tail.Find.Execute("re.Identity.EntityFrameworkCore", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$tail.Font.Color = 255
$tail2 = $d.Content
$tail2.Find.Execute("re.Identity.EntityFrameworkCore", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$tail2.Font.Color = 6299648
